$d = $word.ActiveDocument

# The first paragraph currently reads exactly:
#   "This is a Microsoft word document."
# The commit appends three more runs so the paragraph becomes four runs:
#   "This is a Microsoft word document." + " (" + "Changed main" + ")"
# Each appended run must stay a distinct <w:r> with no run
# properties (matching the target diff), so a plain Range.InsertAfter
# (which merges into the previous identically-formatted run) can't be
# used here. Instead we replace the whole paragraph's XML via
# Range.InsertXML, re-specifying the original run plus the three new
# ones, which keeps them as separate runs while leaving the paragraph's
# own identity (paraId/textId/rsid) untouched.

$para = $d.Paragraphs(1)
$range = $para.Range

if ($range.Text -ne "This is a Microsoft word document.`r") {
    throw "Unexpected first paragraph text: $($range.Text)"
}

$packageXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="5ADF5830" w14:textId="42E3A3E7" w:rsidR="00384372" w:rsidRDefault="00094D0B">
            <w:r><w:t>This is a Microsoft word document.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> (</w:t></w:r>
            <w:r><w:t>Changed main</w:t></w:r>
            <w:r><w:t>)</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$range.InsertXML($packageXml) | Out-Null
